$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Fix the sheet view: scroll back to A1 and select G5 ----
$ws.Range("A1").Select()
$ws.Range("G5").Select()

# ---- New label row 13: "8 TeV" ----
$ws.Range("A13").Value = "8 TeV"

# ---- Row 15: header row (copy of row 1) ----
$ws.Range("A15").Value = "V2"
$ws.Range("B15").Value = "Xi"
$ws.Range("C15").Value = "h-h"
$ws.Range("D15").Value = "fsig"
$ws.Range("F15").Value = "v2"
$ws.Range("H15").Value = "sig"
$ws.Range("K15").Value = "Rel Error"

# ---- Row 16: Obs data row (copy of row 2) ----
$ws.Range("A16").Value = "Obs"
$ws.Range("B16").Value = 0.010265
$ws.Range("C16").Value = 0.0050289999999999996
$ws.Range("D16").Value = 0.95599999999999996

# ---- Row 17: Bkg data row (copy of row 3) ----
$ws.Range("A17").Value = "Bkg"
$ws.Range("B17").Value = 0.0102398

# ---- Row 19: Errors header row (copy of row 5) ----
$ws.Range("A19").Value = "Errors"
$ws.Range("F19").Value = "v2 Errors squared"
$ws.Range("H19").Value = "sig error squared"
$ws.Range("J19").Value = "sig error"

# ---- Row 20: Obs error row (copy of row 6) ----
$ws.Range("A20").Value = "Obs"
$ws.Range("B20").Value = 0.000169605
$ws.Range("C20").Value = 0.00000059501299999999996

# ---- Row 21: Bkg error row (copy of row 7) ----
$ws.Range("A21").Value = "Bkg"
$ws.Range("B21").Value = 0.00028015599999999999

# ---- Merge the C20:C21 pair like C6:C7 (do this before the final
# per-cell formatting pass below, since merging broadcasts the
# top-left cell's style onto the rest of the merged range) ----
$ws.Range("C20:C21").Merge()

# ---- Formulas: set all of them while precedents are still unformatted,
# so this runtime's number-format auto-inheritance doesn't leak onto
# cells that shouldn't receive it (matches rows 2,3,6,7 pattern). ----
$ws.Range("F16").Formula = "=B16/SQRT(C16)"
$ws.Range("H16").Formula = "=(F16-(1-D16)*F17)/(D16)"
$ws.Range("K16").Formula = "=J20/H16"

$ws.Range("F17").Formula = "=B17/SQRT(C16)"

$ws.Range("F20").Formula = "=B20^2/C16 + C20^2*0.25*B16^2/(C16)^3"
$ws.Range("H20").Formula = "=F20/(D16)^2 + F21*((1-D16)/D16)^2"
$ws.Range("J20").Formula = "=SQRT(H20)"

$ws.Range("F21").Formula = "=B21^2/C16 + C20^2*0.25*B17^2/(C16)^3"

$ws.Range("F23").Formula = "=SQRT(F20)"
$ws.Range("F24").Formula = "=SQRT(F21)"

# ---- Now apply number formats / alignment to match the "0.00E+00"
# (style index 1) look used throughout, and the special paired
# numFmt+center / center-only look used by the C6:C7 merge ----
$ws.Range("B16").NumberFormat = "0.00E+00"
$ws.Range("C16").NumberFormat = "0.00E+00"
$ws.Range("F16").NumberFormat = "0.00E+00"

$ws.Range("B17").NumberFormat = "0.00E+00"
$ws.Range("F17").NumberFormat = "0.00E+00"

$ws.Range("B20").NumberFormat = "0.00E+00"
$ws.Range("C20").NumberFormat = "0.00E+00"
$ws.Range("C20").HorizontalAlignment = -4108
$ws.Range("F20").NumberFormat = "0.00E+00"

$ws.Range("B21").NumberFormat = "0.00E+00"
$ws.Range("C21").HorizontalAlignment = -4108
$ws.Range("F21").NumberFormat = "0.00E+00"

$wb.Save()
